$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage for numeric-looking Price values so they stay strings
# (matches source workbook, which stores Price as text) and preserve exact formatting.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Update Price (D) / Volume(1h) (E) cells with refreshed figures
$ws.Range("D2").Value = "64.058.16"
$ws.Range("E2").Value = "  +1.07%  "
$ws.Range("D3").Value = "2.646.42"
$ws.Range("E3").Value = "  -0.44%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "607.35"
$ws.Range("E5").Value = "  -0.44%  "
$ws.Range("D6").Value = "148.64"
$ws.Range("E6").Value = "  +3.35%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("E8").Value = "  +0.73%  "
$ws.Range("E9").Value = "  +1.68%  "
$ws.Range("E10").Value = "  +7.07%  "
$ws.Range("E11").Value = "  -0.14%  "
$ws.Range("E12").Value = "  -0.86%  "
$ws.Range("E13").Value = "  +0.79%  "
$ws.Range("D14").Value = "3.124.73"
$ws.Range("E14").Value = "  -0.27%  "
$ws.Range("D15").Value = "63.928.56"
$ws.Range("E15").Value = "  +1.10%  "
$ws.Range("E16").Value = "  +1.94%  "
$ws.Range("D17").Value = "2.650.22"
$ws.Range("E17").Value = "  -0.62%  "
$ws.Range("D18").Value = "11.95"
$ws.Range("E18").Value = "  +4.41%  "
$ws.Range("E19").Value = "  +3.69%  "
$ws.Range("D20").Value = "346.64"
$ws.Range("E20").Value = "  +1.37%  "
$ws.Range("E21").Value = "  +0.65%  "
$ws.Range("E22").Value = "  +0.08%  "
$ws.Range("E23").Value = "  -0.80%  "
$ws.Range("D24").Value = "66.32"
$ws.Range("E24").Value = "  -0.95%  "
$ws.Range("D25").Value = "1.68"
$ws.Range("E25").Value = "  +8.31%  "
$ws.Range("E26").Value = "  +4.44%  "
$ws.Range("D27").Value = "9.34"
$ws.Range("E27").Value = "  +7.92%  "
$ws.Range("D28").Value = "556.61"
$ws.Range("E28").Value = "  +1.75%  "
$ws.Range("E29").Value = "  +4.10%  "
$ws.Range("E30").Value = "  -1.29%  "
$ws.Range("E31").Value = "  +0.28%  "
$ws.Range("E32").Value = "  +1.08%  "
$ws.Range("E33").Value = "  +5.15%  "
$ws.Range("E34").Value = "  -0.95%  "
$ws.Range("D35").Value = "5.34"
$ws.Range("E35").Value = "  +3.82%  "
$ws.Range("D36").Value = "168.41"
$ws.Range("E36").Value = "  -2.66%  "
$ws.Range("E37").Value = "  +0.17%  "
$ws.Range("E38").Value = "  -0.01%  "
$ws.Range("E39").Value = "  +4.83%  "
$ws.Range("D40").Value = "19.30"
$ws.Range("E40").Value = "  +0.84%  "
$ws.Range("E41").Value = "  +0.06%  "
$ws.Range("D42").Value = "167.26"
$ws.Range("E42").Value = "  -4.17%  "
$ws.Range("D43").Value = "40.30"
$ws.Range("E43").Value = "  +0.50%  "
$ws.Range("D44").Value = "3.83"
$ws.Range("E44").Value = "  +2.20%  "
$ws.Range("D45").Value = "22.03"
$ws.Range("E45").Value = "  -0.82%  "
$ws.Range("E46").Value = "  -0.26%  "
$ws.Range("E47").Value = "  -0.65%  "
$ws.Range("D48").Value = "0.0246"
$ws.Range("E48").Value = "  +2.55%  "
$ws.Range("E49").Value = "  +13.37%  "
$ws.Range("D50").Value = "0.0962"
$ws.Range("E50").Value = "  +0.03%  "
$ws.Range("D51").Value = "19.03"
$ws.Range("E51").Value = "  +1.56%  "
